function Find-ParaByText($doc, $text) {
    $r = $doc.Content
    $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Paragraph not found: $text"
    }
    return $r.Paragraphs(1)
}

# Inserts a sequence of (style, text) paragraphs right after $afterPara.
# Returns the last paragraph inserted (so callers can keep chaining if needed).
function Insert-Paragraphs($afterPara, $items) {
    $cur = $afterPara
    foreach ($item in $items) {
        $cur.Range.InsertParagraphAfter()
        $cur = $cur.Next()
        if ($item.style) {
            $cur.Style = $item.style
        }
        $cur.Range.Text = $item.text
    }
    return $cur
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the contact-info paragraph (phone/email/site/linkedin/location)
#    that sat directly under the centered name heading.
# ---------------------------------------------------------------------------
$contact = Find-ParaByText $d "+1 (512) 555-0123"
$contact.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Professional summary rewrite.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Experienced data scientist and software engineer with 15+ years of expertise in geospatial analysis, demographic research, and political data. Proven track record of building scalable systems, conducting complex analyses, and delivering actionable insights for campaigns, organizations, and government agencies.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Senior data scientist and software engineer specializing in geospatial machine learning and large-scale demographic analysis. Developed algorithms that improved demographic classification accuracy from 23% to 64%, processed data across 178,000+ precincts, and built platforms serving thousands of analysts nationwide.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Core competencies line is cleared out (heading remains, detail blanked).
# ---------------------------------------------------------------------------
$coreLine = Find-ParaByText $d "CODE"
$coreLine.Range.Text = ""

# ---------------------------------------------------------------------------
# 4. Siege Analytics block: update title dates, subtitle, and all 3 bullets.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Partner - Siege Analytics (Austin, TX) | 2020 - Present",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Partner - Siege Analytics (Austin, TX) | 2005 - Present",
    2) | Out-Null

$d.Content.Find.Execute(
    "Data Science & Political Analytics",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Data, Technology and Strategy Consulting",
    2) | Out-Null

$d.Content.Find.Execute(
    "Uncovered decades of demographic miscoding in voter files, discovering 2.7M previously mischaracterized Democratic voters",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%",
    2) | Out-Null

$d.Content.Find.Execute(
    "Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration",
    2) | Out-Null

$d.Content.Find.Execute(
    "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct redistricting analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Replace the remainder of PROFESSIONAL EXPERIENCE (everything from the
#    old "Senior Data Scientist - Lake Research Partners" entry through the
#    old "Data Analyst - The Feldman Group" entry) with the new 7-job block.
# ---------------------------------------------------------------------------
$blockStart = Find-ParaByText $d "Senior Data Scientist - Lake Research Partners"
$blockEndAnchor = Find-ParaByText $d "Trained staff on PHP/MySQL for data analysis and reporting systems"

$delStart = $blockStart.Range.Start
$delEnd = $blockEndAnchor.Range.End
$d.Range($delStart, $delEnd).Delete()

# Anchor for insertion = the last bullet of the Siege Analytics / Partner block.
$expAnchor = Find-ParaByText $d "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins"

$newJobs = @(
    @{style="Heading 3"; text="Data Products Manager - Helm/Murmuration (Austin, TX) | June 2021 - May 2023"},
    @{style="Normal"; text="Civic Graph & Civic Pulse Director"},
    @{style="Normal"; text="• Conceived, architected and built Civic Graph multi-tenant data warehouse processing government data from Census, Bureau of Labor Statistics, National Council of Educational Statistics"},
    @{style="Normal"; text="• Built multi-dimensional data warehouse measuring socio-economic changes in America at every level across attitudinal, behavioral, demographic, economic and geographical dimensions"},
    @{style="Normal"; text="• Managed engineering teams of 7-11 professionals while setting technical direction for data architecture"},

    @{style="Heading 3"; text="Analytics Supervisor - GSD&M (Austin, TX) | November 2019 - June 2020"},
    @{style="Normal"; text="Big Data Engineering Transformation"},
    @{style="Normal"; text="• Transformed small data team into big data engineering team, scaling from laptop datasets to Hadoop Clusters and Hive on AWS"},
    @{style="Normal"; text="• Managed accounts including United States Air Force, Southwest Airlines/Chase and Indeed"},
    @{style="Normal"; text="• Rewrote mission and offerings of department and drafted integration plan with strategy team"},

    @{style="Heading 3"; text="Software Engineer - Mautinoa Technologies (Austin, TX) | August 2016 - February 2018"},
    @{style="Normal"; text="SimCrisis Product Owner/Engineer"},
    @{style="Normal"; text="• Conceived, architected and engineered econometric simulation software for humanitarian crises intervention measurement"},
    @{style="Normal"; text="• Built SimCrisis GeoDjango web application using multi-agent modeling to create econometric simulations of crisis economies"},
    @{style="Normal"; text="• Designed modular application accepting rules extensions for ethnic strife, different crises/disasters, supply failures"},

    @{style="Heading 3"; text="Senior Analyst - Myers Research (Austin, TX) | August 2012 - February 2014"},
    @{style="Normal"; text="RACSO Product Owner/Engineer"},
    @{style="Normal"; text="• Designed comprehensive survey instruments for specialized voting segments and niche markets"},
    @{style="Normal"; text="• Co-developed RACSO web application managing all aspects of survey operations from instrument design to data analysis"},
    @{style="Normal"; text="• Wrote RFP and analyzed bids from 1,200 vendors for research platform development"},

    @{style="Heading 3"; text="Research Director - PCCC (Washington, DC) | 2010 - 2012"},
    @{style="Normal"; text="Political Research & Data Analysis (FLEEM System)"},
    @{style="Normal"; text="• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys"},
    @{style="Normal"; text="• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren"},
    @{style="Normal"; text="• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver"},

    @{style="Heading 3"; text="Software Engineer - Salsa Labs (Washington, DC) | January 2011 - August 2011"},
    @{style="Normal"; text="Geospatial CRM Development"},
    @{style="Normal"; text="• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands simultaneously"},
    @{style="Normal"; text="• Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers"},
    @{style="Normal"; text="• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill"},

    @{style="Heading 3"; text="Programmer - Lake Research Partners (Washington, DC) | April 2008 - December 2008"},
    @{style="Normal"; text="Political Research & Analytics"},
    @{style="Normal"; text="• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party"},
    @{style="Normal"; text="• Harmonized data from 20+ polling firms with incompatible methodologies and encoding systems"},
    @{style="Normal"; text="• Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+"}
)

Insert-Paragraphs $expAnchor $newJobs | Out-Null

Write-Output "stage3-done"
